$wb = $excel.ActiveWorkbook

# Shared date/time number format used by column A throughout the workbook
$dtFormat = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------
# Sheet "ROW50-FE-LIFTER": append row 37
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ROW50-FE-LIFTER")
$r = 37
$ws.Cells.Item($r, 1).NumberFormat = $dtFormat
$ws.Cells.Item($r, 1).Value = 45742.66588601852
$ws.Cells.Item($r, 2).Value = "0x01,0x90"
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x72"
$ws.Cells.Item($r, 5).Value = "0xe"
$ws.Cells.Item($r, 6).Value = 400
$ws.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item($r, 8).Value = 370
$ws.Cells.Item($r, 9).Value = 14

# ---------------------------------------------------------------------
# Sheet "ROW50-MID-LIFTER": append row 39
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ROW50-MID-LIFTER")
$r = 39
$ws.Cells.Item($r, 1).NumberFormat = $dtFormat
$ws.Cells.Item($r, 1).Value = 45742.63417824074
$ws.Cells.Item($r, 2).Value = "0x01,0x90 "
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x76"
$ws.Cells.Item($r, 5).Value = "0x19"
$ws.Cells.Item($r, 6).Value = 400
# This value overflows double precision if parsed as a number, so it must
# stay a text value. Copy it from the identical value already stored as
# text in the row above (row 38) rather than assigning a numeric-looking
# string literal, which Excel would auto-convert to a (lossy) Number.
$ws.Cells.Item(38, 7).Copy($ws.Cells.Item($r, 7))
$ws.Cells.Item($r, 8).Value = 374
$ws.Cells.Item($r, 9).Value = 25

# ---------------------------------------------------------------------
# Sheet "ROW11-FE-LIFTER": append row 37
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ROW11-FE-LIFTER")
$r = 37
$ws.Cells.Item($r, 1).NumberFormat = $dtFormat
$ws.Cells.Item($r, 1).Value = 45742.68198682871
$ws.Cells.Item($r, 2).Value = "0x01,0x90"
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x72"
$ws.Cells.Item($r, 5).Value = "0x14"
$ws.Cells.Item($r, 6).Value = 400
$ws.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item($r, 8).Value = 370
$ws.Cells.Item($r, 9).Value = 20

# ---------------------------------------------------------------------
# Sheet "ROW11-MID-LIFTER": append row 37
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ROW11-MID-LIFTER")
$r = 37
$ws.Cells.Item($r, 1).NumberFormat = $dtFormat
$ws.Cells.Item($r, 1).Value = 45742.83058734953
$ws.Cells.Item($r, 2).Value = "0x01,0x90"
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x76"
$ws.Cells.Item($r, 5).Value = "0x19"
$ws.Cells.Item($r, 6).Value = 400
$ws.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item($r, 8).Value = 374
$ws.Cells.Item($r, 9).Value = 25

Write-Host "Appended new log rows to all four sheets"
